$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.242767026
$ws.Range("E4").Value = 1.089141008
$ws.Range("E5").Value = 1.202749768
$ws.Range("E6").Value = 1.16005959
$ws.Range("E7").Value = 1.453985293
$ws.Range("E8").Value = 1.448976496
$ws.Range("E9").Value = 1.47877001
$ws.Range("E10").Value = 1.434887836
$ws.Range("E11").Value = 1.469273649
$ws.Range("E12").Value = 1.267817051
$ws.Range("E13").Value = 1.009392797
$ws.Range("E14").Value = 1.368467369
$ws.Range("E16").Value = 1.293594346
